$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.182.24'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.826.07'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9982'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.20'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6091'
$ws.Range('E6').Value = '  -3.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9998'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07105'
$ws.Range('E8').Value = '  -4.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2813'
$ws.Range('E9').Value = '  -2.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.65'
$ws.Range('E10').Value = '  -5.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07657'
$ws.Range('E11').Value = '  -0.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.780.21'
$ws.Range('E12').Value = '  -3.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.823'
$ws.Range('E13').Value = '  -2.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.00001007'
$ws.Range('E14').Value = '  -1.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6323'
$ws.Range('E15').Value = '  -6.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.069.52'
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.07'
$ws.Range('E17').Value = '  -3.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.872'
$ws.Range('E18').Value = '  -6.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.157.73'
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '227.67'
$ws.Range('E20').Value = '  -0.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.79'
$ws.Range('E21').Value = '  -4.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9995'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.012'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9994'
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.79'
$ws.Range('E25').Value = '  -2.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.052'
$ws.Range('E26').Value = '  -5.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1310'
$ws.Range('E27').Value = '  -2.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.60'
$ws.Range('E28').Value = '  -4.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.487'
$ws.Range('E29').Value = '  +1.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06448'
$ws.Range('E30').Value = '  -8.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.454'
$ws.Range('E31').Value = '  -1.59%  '
$ws.Range('E32').Value = '  -5.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.797'
$ws.Range('E33').Value = '  -6.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.126'
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.749'
$ws.Range('E35').Value = '  -4.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6505'
$ws.Range('E36').Value = '  -6.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.544'
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.219.77'
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.748'
$ws.Range('E39').Value = '  -2.35%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01744'
$ws.Range('E40').Value = '  -5.41%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.571'
$ws.Range('E41').Value = '  -3.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9286'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9988'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.17'
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.974.10'
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.04'
$ws.Range('E46').Value = '  -3.46%  '
$ws.Range('E47').Value = '  -0.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.623'
$ws.Range('E48').Value = '  -4.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.587'
$ws.Range('E49').Value = '  -3.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4559'
$ws.Range('E50').Value = '  -0.70%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05519'
$ws.Range('E51').Value = '  -2.71%  '

Write-Host "Applied 107 cell updates"
